$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection(1)
$s1.Formula = '=SERIES("Start Date",Sheet1!$A$2:$A$18,Sheet1!$B$2:$B$18,1)'
Write-Output $s1.Values.Count
